$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.54194700717926
$ws.Range("B1").Value = 2.218854665756226
$ws.Range("C1").Value = -1
$ws.Range("D1").Value = 1.369795441627502
$ws.Range("E1").Value = 0.6508402824401855
